$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6925.25
$ws.Range("J64").Value = 7334.3335
$ws.Range("L64").Value = 7334.3335
$ws.Range("N64").Value = -7830.3335
$ws.Range("H67").Value = 6925.25
$ws.Range("J67").Value = 7334.3335
$ws.Range("L67").Value = 7334.3335
$ws.Range("N67").Value = -9050.333500000001
$ws.Range("H70").Value = 1962.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 1962.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 5887.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -6427.5
$ws.Range("H73").Value = 1962.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 1962.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 5887.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -7759.5
$ws.Range("H138").Value = 3577.7654
$ws.Range("I138").Value = 2673.4
$ws.Range("J138").Value = 3741.2048
$ws.Range("K138").Value = 8020.200000000001
$ws.Range("L138").Value = 11223.6144
$ws.Range("M138").Value = -2880.200000000001
$ws.Range("N138").Value = -21503.6144

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 30044
$ws.Range("I54").Value = 30044
$ws.Range("K54").Value = 30044
$ws.Range("M54").Value = -29275
$ws.Range("H63").Value = 5575
$ws.Range("I63").Value = 2500
$ws.Range("J63").Value = 6600
$ws.Range("K63").Value = 2500
$ws.Range("L63").Value = 6600
$ws.Range("M63").Value = -1814
$ws.Range("N63").Value = -7972
$ws.Range("H66").Value = 5575
$ws.Range("I66").Value = 2500
$ws.Range("J66").Value = 6600
$ws.Range("K66").Value = 12500
$ws.Range("L66").Value = 33000
$ws.Range("M66").Value = -9068
$ws.Range("N66").Value = -39864
$ws.Range("H81").Value = 121189.8
$ws.Range("J81").Value = 121987.25
$ws.Range("L81").Value = 121987.25
$ws.Range("N81").Value = -123983.25
$ws.Range("H84").Value = 121189.8
$ws.Range("J84").Value = 121987.25
$ws.Range("L84").Value = 365961.75
$ws.Range("N84").Value = -375945.75
$ws.Range("H92").Value = 127999
$ws.Range("J92").Value = 127999
$ws.Range("L92").Value = 127999
$ws.Range("N92").Value = -132991
$ws.Range("H122").Value = 3642.4783
$ws.Range("I122").Value = 2872.6875
$ws.Range("J122").Value = 5402
$ws.Range("K122").Value = 8618.0625
$ws.Range("L122").Value = 16206
$ws.Range("M122").Value = -6168.0625
$ws.Range("N122").Value = -21106

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1740.5769
$ws.Range("I20").Value = 1871.8235
$ws.Range("K20").Value = 1871.8235
$ws.Range("M20").Value = -1624.8235

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 166667420
$ws.Range("I16").Value = 166667420
$ws.Range("K16").Value = 166667420
$ws.Range("M16").Value = -166667133
$ws.Range("H31").Value = 1973.5103
$ws.Range("I31").Value = 1382.0952
$ws.Range("K31").Value = 1382.0952
$ws.Range("M31").Value = -1087.0952
$ws.Range("H34").Value = 1973.5103
$ws.Range("I34").Value = 1382.0952
$ws.Range("K34").Value = 1382.0952
$ws.Range("M34").Value = -1180.0952
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 25000
$ws.Range("K42").Value = 25000
$ws.Range("M42").Value = -24407
$ws.Range("H113").Value = 166667420
$ws.Range("I113").Value = 166667420
$ws.Range("K113").Value = 166667420
$ws.Range("M113").Value = -166665250
$ws.Range("H132").Value = 3868.0557
$ws.Range("I132").Value = 3687.5334
$ws.Range("J132").Value = 4770.6665
$ws.Range("K132").Value = 11062.6002
$ws.Range("L132").Value = 14311.9995
$ws.Range("M132").Value = -8532.600199999999
$ws.Range("N132").Value = -19371.9995
$ws.Range("H134").Value = 2118.5938
$ws.Range("I134").Value = 2118.5938
$ws.Range("K134").Value = 6355.7814
$ws.Range("M134").Value = -3820.7814

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 950.75
$ws.Range("J34").Value = 1401
$ws.Range("L34").Value = 4203
$ws.Range("N34").Value = -4371
$ws.Range("H68").Value = 2506.8572
$ws.Range("I68").Value = 1584.7142
$ws.Range("J68").Value = 3429
$ws.Range("K68").Value = 4754.142599999999
$ws.Range("L68").Value = 10287
$ws.Range("M68").Value = -3943.142599999999
$ws.Range("N68").Value = -11909
$ws.Range("H71").Value = 2506.8572
$ws.Range("I71").Value = 1584.7142
$ws.Range("J71").Value = 3429
$ws.Range("K71").Value = 14262.4278
$ws.Range("L71").Value = 30861
$ws.Range("M71").Value = -10206.4278
$ws.Range("N71").Value = -38973
$ws.Range("H88").Value = 3983.1936
$ws.Range("J88").Value = 3983.1936
$ws.Range("L88").Value = 11949.5808
$ws.Range("N88").Value = -12805.5808
$ws.Range("H91").Value = 3983.1936
$ws.Range("J91").Value = 3983.1936
$ws.Range("L91").Value = 11949.5808
$ws.Range("N91").Value = -14913.5808

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 26669.385
$ws.Range("I113").Value = 3442.2
$ws.Range("J113").Value = 41186.375
$ws.Range("K113").Value = 3442.2
$ws.Range("L113").Value = 41186.375
$ws.Range("M113").Value = -1272.2
$ws.Range("N113").Value = -45526.375
$ws.Range("H122").Value = 1867.75
$ws.Range("I122").Value = 1728.6666
$ws.Range("J122").Value = 2285
$ws.Range("K122").Value = 5185.9998
$ws.Range("L122").Value = 6855
$ws.Range("M122").Value = -2735.9998
$ws.Range("N122").Value = -11755
$ws.Range("H139").Value = 457133.2
$ws.Range("J139").Value = 457133.2
$ws.Range("L139").Value = 457133.2
$ws.Range("N139").Value = -467413.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3920
$ws.Range("I61").Value = 1603.8182
$ws.Range("K61").Value = 1603.8182
$ws.Range("M61").Value = -1401.8182
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 3920
$ws.Range("I113").Value = 1603.8182
$ws.Range("K113").Value = 1603.8182
$ws.Range("M113").Value = 566.1818000000001
$ws.Range("H132").Value = 8845.200000000001
$ws.Range("I132").Value = 8845.200000000001
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 26535.6
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -24005.6
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 129999
$ws.Range("J133").Value = 129999
$ws.Range("L133").Value = 129999
$ws.Range("N133").Value = -135059
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 61732.832
$ws.Range("I136").Value = 90124.25
$ws.Range("K136").Value = 270372.75
$ws.Range("M136").Value = -267822.75
$ws.Range("H141").Value = 130500
$ws.Range("J141").Value = 130500
$ws.Range("L141").Value = 130500
$ws.Range("N141").Value = -140860

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 56580.855
$ws.Range("J54").Value = 74999.2
$ws.Range("L54").Value = 74999.2
$ws.Range("N54").Value = -76039.2
$ws.Range("H62").Value = 5803.75
$ws.Range("I62").Value = 5706.4287
$ws.Range("K62").Value = 5706.4287
$ws.Range("M62").Value = -5082.4287
$ws.Range("H65").Value = 5803.75
$ws.Range("I65").Value = 5706.4287
$ws.Range("K65").Value = 28532.1435
$ws.Range("M65").Value = -25412.1435
$ws.Range("H136").Value = 38908.105
$ws.Range("J136").Value = 80292.38
$ws.Range("L136").Value = 240877.14
$ws.Range("N136").Value = -245977.14
